$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 98.04528151073291
$ws.Range("C7").Value = 98.03921256460907
$ws.Range("D7").Value = 98.07170447330606
$ws.Range("E7").Value = 98.01109088892491

$ws.Range("B8").Value = 97.65581094026977
$ws.Range("C8").Value = 97.53038799035613
$ws.Range("D8").Value = 97.66796705025925
$ws.Range("E8").Value = 97.58771602536946

$ws.Range("B9").Value = 96.24429554855158
$ws.Range("C9").Value = 96.3317856062
$ws.Range("D9").Value = 96.27242821750166
$ws.Range("E9").Value = 96.28221152061533
